$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so that numeric-looking
# strings such as "309.02" are written back verbatim instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.846.29"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.838.08"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "309.02"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4704"
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("D8").Value = "0.3653"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "0.07143"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "0.9206"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").Value = "1.980.13"
$ws.Range("E11").Value = "  +9.02%  "
$ws.Range("D12").Value = "19.54"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "0.07654"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "5.282"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "6.389"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "88.09"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "0.000008629"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "26.881.50"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "14.44"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").Value = "5.005"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "1.919"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "151.51"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "18.20"
$ws.Range("D27").Value = "2.006"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").Value = "114.08"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").Value = "0.08816"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").Value = "1.175"
$ws.Range("E32").Value = "  +5.96%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "2.748"
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").Value = "4.469"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "0.01942"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").Value = "0.05207"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "0.5193"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("D41").Value = "6.953"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "0.1511"
$ws.Range("D43").Value = "8.146"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "10.44"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D45").Value = "0.4695"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D47").Value = "101.42"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("D48").Value = "1.592"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "64.94"
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("D50").Value = "0.06034"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "0.8858"
$ws.Range("E51").Value = "  +4.75%  "

# Restore the default (Normal) style so the cells keep their original,
# unformatted appearance once the text values have been written.
$ws.Range("D2:D51").Style = "Normal"
